$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new record row (row 2) under the header row.
$ws.Range("A2").Value = "MCH141"
$ws.Range("C2").Value = "LIST OF NAMES/PRISON NUMBERS, WOMEN POLITICAL PRISONERS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Apply the record-row formatting (10pt Calibri, automatic/theme text color) to every
# cell in the row, including the still-empty D2/H2 cells (date_s / file_path columns).
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1

$ws.Range("C2:H2").Font.Name = "Calibri"
$ws.Range("C2:H2").Font.Size = 10
$ws.Range("C2:H2").Font.ThemeColor = 1

# Freeze the header row and select the new record row, matching the saved view state.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A2:H2").Select()
